$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1563
$ws1.Range("F7").Value = 214
$ws1.Range("F10").Value = 1137
$ws1.Range("F11").Value = 667
$ws1.Range("F12").Value = 448
$ws1.Range("F13").Value = 738
$ws1.Range("F14").Value = 71
$ws1.Range("F15").Value = 213
$ws1.Range("F16").Value = 181
$ws1.Range("F17").Value = 230
$ws1.Range("F18").Value = 151
$ws1.Range("F19").Value = 278
$ws1.Range("F20").Value = 1468
$ws1.Range("F22").Value = 66
$ws1.Range("F25").Value = 2208
$ws1.Range("F27").Value = 774
$ws1.Range("F29").Value = 72
$ws1.Range("F30").Value = 44

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 74
$ws2.Range("F15").Value = 436

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value = 578
$ws3.Range("F9").Value = 443

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 1563
$ws4.Range("F12").Value = 578
$ws4.Range("F14").Value = 443
$ws4.Range("F16").Value = 214
$ws4.Range("F19").Value = 1137
$ws4.Range("F20").Value = 667
$ws4.Range("F21").Value = 448
$ws4.Range("F24").Value = 738
$ws4.Range("F25").Value = 71
$ws4.Range("F26").Value = 213
$ws4.Range("F27").Value = 74
$ws4.Range("F30").Value = 181
$ws4.Range("F31").Value = 230
$ws4.Range("F32").Value = 151
$ws4.Range("F33").Value = 278
$ws4.Range("F35").Value = 1468
$ws4.Range("F38").Value = 66
$ws4.Range("F41").Value = 2208
$ws4.Range("F44").Value = 774
$ws4.Range("F47").Value = 72
$ws4.Range("F48").Value = 44
